$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Formula = "=2+3"
$ws.Range("C7").Formula = "=3+6"
$ws.Range("C8").Formula = "=1=8"

$null = $ws.Range("C9").Select()
